$wb = $excel.ActiveWorkbook

# Sheet "2025" (first sheet)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2778.902526399997
$ws.Range("E2").Value = 290927.2506141524
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312751
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 71977.22211760026
$ws.Range("O2").Value = 68708.80120585448

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52597.49431690662
$ws.Range("E2").Value = 269323.7120331937
$ws.Range("I2").Value = 221949.8854910079
$ws.Range("L2").Value = 225013.3978349316
$ws.Range("M2").Value = 105604.6794510125
$ws.Range("N2").Value = 35977.50527378642
$ws.Range("O2").Value = 25114.09096912013

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22664.87971660737
$ws.Range("B2").Value = 15747.05168594091
$ws.Range("E2").Value = 110396.7489265018
$ws.Range("I2").Value = 188550.6342719396
$ws.Range("M2").Value = 58163.04685719004
$ws.Range("N2").Value = 49759.46314991143
$ws.Range("O2").Value = 47598.28305178237

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("O2").Value = 22386.41312234465
